$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=71; A="10.45.3.7"; B="Srv_Geo_Bk"; C="10.45.3.7"; D="ansible_test"; E="alba5-bk"},
    @{Row=72; A="10.178.3.4"; B="Srv_Geo_Bk"; C="10.178.3.4"; D="ansible_test"; E="albuc-n1"},
    @{Row=73; A="10.98.3.4"; B="Srv_Geo_Bk"; C="10.98.3.4"; D="ansible_test"; E="alrio-n1"},
    @{Row=74; A="10.106.3.4"; B="Srv_Geo_Bk"; C="10.106.3.4"; D="ansible_test"; E="aldie-n1"},
    @{Row=75; A="10.112.3.4"; B="Srv_Geo_Bk"; C="10.112.3.4"; D="ansible_test"; E="alper-n1"},
    @{Row=76; A="10.139.2.3"; B="Srv_Geo_Bk"; C="10.139.2.3"; D="ansible_test"; E="alct2-bk"},
    @{Row=77; A="10.134.3.4"; B="Srv_Geo_Bk"; C="10.134.3.4"; D="ansible_test"; E="alflo-n1"},
    @{Row=78; A="10.187.3.4"; B="Srv_Geo_Bk"; C="10.187.3.4"; D="ansible_test"; E="alca5-n1"},
    @{Row=79; A="10.189.3.4"; B="Srv_Geo_Bk"; C="10.189.3.4"; D="ansible_test"; E="alarm-n1"},
    @{Row=80; A="10.34.3.3"; B="Srv_Geo_Bk"; C="10.34.3.3"; D="ansible_test"; E="alave-n1"},
    @{Row=81; A="10.65.3.4"; B="Srv_Geo_Bk"; C="10.65.3.4"; D="ansible_test"; E="ktman-n1"},
    @{Row=82; A="10.190.3.4"; B="Srv_Geo_Bk"; C="10.190.3.4"; D="ansible_test"; E="ktark-n1"},
    @{Row=83; A="10.231.3.4"; B="Srv_Geo_Bk"; C="10.231.3.4"; D="ansible_test"; E="ktgir-n1"},
    @{Row=84; A="10.92.3.4"; B="Srv_Geo_Bk"; C="10.92.3.4"; D="ansible_test"; E="ktpob-bk"},
    @{Row=85; A="172.16.4.4"; B="Srv_Geo_Bk"; C="172.16.4.4"; D="ansible_test"; E="akven-n1"},
    @{Row=86; A="10.115.3.4"; B="Srv_Geo_Bk"; C="10.115.3.4"; D="ansible_test"; E="kttes-bk"},
    @{Row=87; A="10.91.3.7"; B="Srv_Geo_Bk"; C="10.91.3.7"; D="ansible_test"; E="ktsal-n2"},
    @{Row=88; A="10.147.3.5"; B="Srv_Geo_Bk"; C="10.147.3.5"; D="ansible_test"; E="ktuno-n1"},
    @{Row=89; A="10.127.0.7"; B="Srv_Geo_Bk"; C="10.127.0.7"; D="ansible_test"; E="ktjul-n1"},
    @{Row=90; A="10.203.3.4"; B="Srv_Geo_Bk"; C="10.203.3.4"; D="ansible_test"; E="ktfus-n1"},
    @{Row=91; A="172.16.6.8"; B="Srv_Geo_Bk"; C="172.16.6.8"; D="ansible_test"; E="akbol-bk"},
    @{Row=92; A="172.16.8.3"; B="Srv_Geo_Bk"; C="172.16.8.3"; D="ansible_test"; E="akipi-bk"},
    @{Row=93; A="172.16.9.3"; B="Srv_Geo_Bk"; C="172.16.9.3"; D="ansible_test"; E="akper-n1"},
    @{Row=94; A="172.16.10.6"; B="Srv_Geo_Bk"; C="172.16.10.6"; D="ansible_test"; E="ak170-n1"},
    @{Row=95; A="10.156.3.7"; B="Srv_Geo_Bk"; C="10.156.3.7"; D="ansible_test"; E="alva2-new"},
    @{Row=96; A="10.155.3.4"; B="Srv_Geo_Bk"; C="10.155.3.4"; D="ansible_test"; E="alpie-n1"},
    @{Row=97; A="10.144.3.7"; B="Srv_Geo_Bk"; C="10.144.3.7"; D="ansible_test"; E="aliba-bk"},
    @{Row=98; A="10.149.3.7"; B="Srv_Geo_Bk"; C="10.149.3.7"; D="ansible_test"; E="alnei-bk"},
    @{Row=99; A="10.29.3.5"; B="Srv_Geo_Bk"; C="10.29.3.5"; D="ansible_test"; E="almol-n2"},
    @{Row=100; A="10.86.0.16"; B="Srv_Geo_Bk"; C="10.86.0.16"; D="ansible_test"; E="alkun-n2"},
    @{Row=101; A="10.87.3.8"; B="Srv_Geo_Bk"; C="10.87.3.8"; D="ansible_test"; E="almay-bk"},
    @{Row=102; A="10.108.3.87"; B="Srv_Geo_Bk"; C="10.108.3.87"; D="ansible_test"; E="albel-new"},
    @{Row=103; A="10.49.3.8"; B="Srv_Geo_Bk"; C="10.49.3.8"; D="ansible_test"; E="alnue-bk"},
    @{Row=104; A="10.71.3.5"; B="Srv_Geo_Bk"; C="10.71.3.5"; D="ansible_test"; E="ktb94-n1"},
    @{Row=105; A="10.32.3.8"; B="Srv_Geo_Bk"; C="10.32.3.8"; D="ansible_test"; E="ktchi-bk"},
    @{Row=106; A="10.109.3.8"; B="Srv_Geo_Bk"; C="10.109.3.8"; D="ansible_test"; E="ktame-bk"},
    @{Row=107; A="10.104.3.3"; B="Srv_Geo_Bk"; C="10.104.3.3"; D="ansible_test"; E="ktmay-n1"},
    @{Row=108; A="10.122.3.3"; B="Srv_Geo_Bk"; C="10.122.3.3"; D="ansible_test"; E="kttit-bk"},
    @{Row=109; A="10.129.3.5"; B="Srv_Geo_Bk"; C="10.129.3.5"; D="ansible_test"; E="ktbuc-bk"},
    @{Row=110; A="10.46.3.8"; B="Srv_Geo_Bk"; C="10.46.3.8"; D="ansible_test"; E="ktsoa-bk"},
    @{Row=111; A="10.48.3.8"; B="Srv_Geo_Bk"; C="10.48.3.8"; D="ansible_test"; E="ktcal-bk"},
    @{Row=112; A="10.145.3.7"; B="Srv_Geo_Bk"; C="10.145.3.7"; D="ansible_test"; E="ktmos-bk"},
    @{Row=113; A="10.151.3.3"; B="Srv_Geo_Bk"; C="10.151.3.3"; D="ansible_test"; E="ktvil-bk"},
    @{Row=114; A="10.89.3.3"; B="Srv_Geo_Bk"; C="10.89.3.3"; D="ansible_test"; E="ktsba-n1"},
    @{Row=115; A="10.232.3.6"; B="Srv_Geo_Bk"; C="10.232.3.6"; D="ansible_test"; E="ktbar-n1"},
    @{Row=116; A="172.16.7.15"; B="Srv_Geo_Bk"; C="172.16.7.15"; D="ansible_test"; E="akpas-bk"},
    @{Row=117; A="172.16.90.3"; B="Srv_Geo_Bk"; C="172.16.90.3"; D="ansible_test"; E="aktuq-n1"},
    @{Row=118; A="172.16.91.7"; B="Srv_Geo_Bk"; C="172.16.91.7"; D="ansible_test"; E="akigp-bk"},
    @{Row=119; A="147127110112"; B="Srv_Geo_Bk"; C="147127110112"; D="ansible_test"; E="akb30-n1"},
    @{Row=120; A="172.16.11.7"; B="Srv_Geo_Bk"; C="172.16.11.7"; D="ansible_test"; E="akb68-n1"},
    @{Row=121; A="172.16.3.6"; B="Srv_Geo_Bk"; C="172.16.3.6"; D="ansible_test"; E="akvil-n1"},
    @{Row=122; A="172.16.107.7"; B="Srv_Geo_Bk"; C="172.16.107.7"; D="ansible_test"; E="akbar-bk"},
    @{Row=123; A="172.16.116.80"; B="Srv_Geo_Bk"; C="172.16.116.80"; D="ansible_test"; E="aksin-bk"},
    @{Row=124; A="10.16.77.4"; B="Srv_Geo_Bk"; C="10.16.77.4"; D="ansible_test"; E="akcan-n1"},
    @{Row=125; A="10.173.3.7"; B="Srv_Geo_Bk"; C="10.173.3.7"; D="ansible_test"; E="alapa-bk"},
    @{Row=126; A="10.121.0.6"; B="Srv_Geo_Bk"; C="10.121.0.6"; D="ansible_test"; E="akyop-n1"},
    @{Row=127; A="172.16.134.4"; B="Srv_Geo_Bk"; C="172.16.134.4"; D="ansible_test"; E="akede-n1"},
    @{Row=128; A="10.125.3.8"; B="Srv_Geo_Bk"; C="10.125.3.8"; D="ansible_test"; E="albu2-bk"},
    @{Row=129; A="10.245.3.4"; B="Srv_Geo_Bk"; C="10.245.3.4"; D="ansible_test"; E="ktnqs-n1"},
)

# Source cell carrying the common data style (fontId=3 "Aptos Narrow", xf index 5)
# used for columns A, B, C, E throughout the data rows.
$formatSource = $ws.Range("A60")

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the standard data-row format onto A, B, C, E (column D keeps the
    # workbook default style, matching the rest of the sheet).
    $formatSource.Copy() | Out-Null
    $ws.Range("A$rowNum").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$rowNum").PasteSpecial(-4122) | Out-Null
    $ws.Range("C$rowNum").PasteSpecial(-4122) | Out-Null
    $ws.Range("E$rowNum").PasteSpecial(-4122) | Out-Null

    if ($rowNum -eq 119) {
        # This single record was captured as a literal number (not text) in
        # the source inventory, so Host/IP land as numbers with a
        # "#,##0" number format instead of the usual text style.
        $ws.Range("A$rowNum").Value = 147127110112
        $ws.Range("C$rowNum").Value = 147127110112
        $ws.Range("A$rowNum").NumberFormat = "#,##0"
        $ws.Range("C$rowNum").NumberFormat = "#,##0"
        $ws.Range("A$rowNum").HorizontalAlignment = -4131
        $ws.Range("C$rowNum").HorizontalAlignment = -4131
        $ws.Range("A$rowNum").VerticalAlignment = -4160
        $ws.Range("C$rowNum").VerticalAlignment = -4160
    } else {
        $ws.Range("A$rowNum").Value = $r.A
        $ws.Range("C$rowNum").Value = $r.C
    }

    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
}

# Restore the view to mirror where the author ended up after the paste.
$ws.Range("C131").Select() | Out-Null

Write-Host "Inserted $($rows.Count) inventory rows (71-129)."
